$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$row2 = @(3, 1, 18.741769, 56.225307, 0.2218531826860132, 0.2218531826860132, 3, 1, 9.266076, 27.798228, 0.5506394579555814, 0.5506394579555816, 173.662655928444, 1562.963903355996, 0.1221611162599469, 0.1221611162599469)
$row3 = @(3, 1, 18.741769, 56.225307, 0.2218531826860132, 0.2218531826860132, 3, 1, 2.181666333333333, 6.544999, 0.1296462026888844, 0.1296462026888844, 40.88828645441033, 367.994578089693, 0.02876242268968498, 0.02876242268968498)
$row4 = @(3, 1, 18.741769, 56.225307, 0.2218531826860132, 0.2218531826860132, 3, 1, 5.380103666666667, 16.140311, 0.3197143393555341, 0.3197143393555341, 100.8326601167197, 907.493941050477, 0.07092964373638133, 0.07092964373638135)
$row5 = @(3, 1, 53.77230066666667, 161.316902, 0.6365224138259964, 0.6365224138259964, 3, 1, 9.266076, 27.798228, 0.5506394579555814, 0.5506394579555816, 498.258224672184, 4484.324022049656, 0.3504943569257249, 0.350494356925725)
$row6 = @(3, 1, 53.77230066666667, 161.316902, 0.6365224138259964, 0.6365224138259964, 3, 1, 2.181666333333333, 6.544999, 0.1296462026888844, 0.1296462026888844, 117.3132180303442, 1055.818962273098, 0.0825227138789031, 0.0825227138789031)
$row7 = @(3, 1, 53.77230066666667, 161.316902, 0.6365224138259964, 0.6365224138259964, 3, 1, 5.380103666666667, 16.140311, 0.3197143393555341, 0.3197143393555341, 289.3005519818358, 2603.704967836522, 0.2035053430213683, 0.2035053430213683)
$row8 = @(3, 1, 11.964182, 35.892546, 0.1416244034879904, 0.1416244034879904, 3, 1, 9.266076, 27.798228, 0.5506394579555814, 0.5506394579555816, 110.861019689832, 997.7491772084879, 0.07798398476990961, 0.07798398476990961)
$row9 = @(3, 1, 11.964182, 35.892546, 0.1416244034879904, 0.1416244034879904, 3, 1, 2.181666333333333, 6.544999, 0.1296462026888844, 0.1296462026888844, 26.10185307527266, 234.916677677454, 0.01836106612029636, 0.01836106612029636)
$row10 = @(3, 1, 11.964182, 35.892546, 0.1416244034879904, 0.1416244034879904, 3, 1, 5.380103666666667, 16.140311, 0.3197143393555341, 0.3197143393555341, 64.36853944686733, 579.316855021806, 0.04527935259778446, 0.04527935259778446)

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10)
$data = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $rows[$r]
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Range($cols[$c] + $rowNum).Value = $rowVals[$c]
    }
}
